$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($i = 0; $i -le 19; $i++) {
    $row = 2 + $i
    $num = "{0:D2}" -f $i
    $ws.Range("A$row").Value = "sequences/278857_learning_sequence_$num.csv"
}
